$d = $word.ActiveDocument

$d.Content.Find.Execute("[[1]]", $true, $false, $false, $false, $false,
                         $true, 1, $false, "`${laudo}", 2)

$d.Content.Find.Execute("[[2]]", $true, $false, $false, $false, $false,
                         $true, 1, $false, "`${analise}", 2)

$d.Content.Find.Execute("[[3]]", $true, $false, $false, $false, $false,
                         $true, 1, $false, "`${consideracoesTecnicas}", 2)
